$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new column K ("2022") to the right of the existing data (which ran
# through column J / year 2021). For every existing row we replicate the
# formatting already used by the row (picking a source cell whose style
# matches what the new K cell needs) and then write in the new value.
# ---------------------------------------------------------------------------

function Copy-Format([string]$src, [string]$dst) {
    $ws.Range($src).Copy() | Out-Null
    $ws.Range($dst).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# Row 3 - thin bottom-border separator row, cell stays empty but keeps style
Copy-Format "J3" "K3"

# Row 4 - header row (years); K4 gets the new year value 2022
Copy-Format "D4" "K4"
$ws.Range("K4").Value2 = 2022

# Row 5 - totals row
Copy-Format "D5" "K5"
$ws.Range("K5").Value2 = 6512.5

# Row 6
Copy-Format "D6" "K6"
$ws.Range("K6").Value2 = 9.2524142373849365

# Row 8
Copy-Format "D8" "K8"
$ws.Range("K8").Value2 = 9.1726448995762642

# Row 9
Copy-Format "D9" "K9"
$ws.Range("K9").Value2 = 9.2988242598562199

# Row 11
Copy-Format "D11" "K11"
$ws.Range("K11").Value2 = 9.3929513987987647

# Row 12
Copy-Format "D12" "K12"
$ws.Range("K12").Value2 = 9.0963110935638873

# Row 14
Copy-Format "D14" "K14"
$ws.Range("K14").Value2 = 12.819848845068858

# Row 15
Copy-Format "D15" "K15"
$ws.Range("K15").Value2 = 7.8843953890120773

# Row 16
Copy-Format "I16" "K16"
$ws.Range("K16").Value2 = 7.6890997954400655

# Row 17
Copy-Format "D17" "K17"
$ws.Range("K17").Value2 = 4.859109337853182

# Row 19
Copy-Format "D19" "K19"
$ws.Range("K19").Value2 = 12.268952512062626

# Row 20
Copy-Format "D20" "K20"
$ws.Range("K20").Value2 = 8.8432516850244731

# Row 21
Copy-Format "D21" "K21"
$ws.Range("K21").Value2 = 12.356872582336921

# Row 22
Copy-Format "D22" "K22"
$ws.Range("K22").Value2 = 25.295368484771757

# Row 23
Copy-Format "D23" "K23"
$ws.Range("K23").Value2 = 4.2612456375718351

# Row 24
Copy-Format "D24" "K24"
$ws.Range("K24").Value2 = 14.933279226285201

# Row 25
Copy-Format "D25" "K25"
$ws.Range("K25").Value2 = 9.0993456624506877

# Row 26
Copy-Format "D26" "K26"
$ws.Range("K26").Value2 = 6.7003522302183303

# Row 27
Copy-Format "D27" "K27"
$ws.Range("K27").Value2 = 12.078370902890091

# ---------------------------------------------------------------------------
# Update the selection to match the new active cell left after the edits.
# ---------------------------------------------------------------------------
$ws.Range("L3").Select() | Out-Null
